$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the underline formatting from D6 (previously a distinct style s="3",
# now matches the plain centered style s="1" used by D2:D5)
$ws.Range("D6").Font.Underline = -4142  # xlUnderlineStyleNone

# Clear the old stray cell M19 (and its formatting) so it no longer
# contributes to the used range / dimension
$ws.Range("M19").Clear()

# Add a new data row (row 7) with the new position, reusing the formatting
# already used by the other data rows (A:D and F centered, E/H default)
$ws.Range("A2:D2").Copy()
$ws.Range("A7:D7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F2").Copy()
$ws.Range("F7").PasteSpecial(-4122)     # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A7").Value = "SEMAPA"
$ws.Range("B7").Value = "SEM.LS"
$ws.Range("C7").Value = "SEM.LS"
$ws.Range("D7").Value = "EUR"
$ws.Range("E7").Value = "STK"
$ws.Range("F7").Value = 1000
$ws.Range("H7").Value = 1

# Update the active selection like the author's session
$ws.Range("K6").Select()
